$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.252.40"
$ws.Range("E2").Value = "  +3.34%  "
$ws.Range("D3").Value = "2.995.06"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.83"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.95"
$ws.Range("E6").Value = "  +5.00%  "
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "2.981.88"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.14"
$ws.Range("E11").Value = "  +7.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +4.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.62"
$ws.Range("E14").Value = "  +3.43%  "
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "3.492.71"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("E17").Value = "  +6.68%  "
$ws.Range("D18").Value = "2.994.55"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").Value = "59.260.45"
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.07"
$ws.Range("E20").Value = "  +3.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").Value = "  +5.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.27"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.47"
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.18"
$ws.Range("E28").Value = "  +10.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.54"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.80"
$ws.Range("E30").Value = "  +4.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.68"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0992"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.993"
$ws.Range("E34").Value = "  +6.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.91"
$ws.Range("E35").Value = "  +5.18%  "
$ws.Range("D36").Value = "0.0₃0758"
$ws.Range("E36").Value = "  +11.78%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.98"
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.64"
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("E40").Value = "  +6.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "399.41"
$ws.Range("E41").Value = "  +6.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0351"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "2.750.37"
$ws.Range("E43").Value = "  +4.38%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.251"
$ws.Range("E45").Value = "  +4.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.48"
$ws.Range("E46").Value = "  +25.82%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.25"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.110"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.34"
$ws.Range("E51").Value = "  +0.23%  "

foreach ($addr in @("D5","D6","D10","D11","D12","D13","D14","D20","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D33","D34","D35","D38","D39","D40","D41","D42","D45","D46","D48","D49","D51")) {
    $ws.Range($addr).Style = "Normal"
}
